# Update the "Elements" timing results (column C) on each sort-algorithm
# sheet to reflect the latest benchmark run recorded in the commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Bubble Sort")
$ws.Range("C2").Value = 6
$ws.Range("C3").Value = 149
$ws.Range("C4").Value = 13191
$ws.Range("C6").Value = 121
$ws.Range("C7").Value = 13209
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 120
$ws.Range("C10").Value = 13205
$ws.Range("C12").Value = 121
$ws.Range("C13").Value = 13173
$ws.Range("C15").Value = 121
$ws.Range("C16").Value = 13270
$ws.Range("C18").Value = 122
$ws.Range("C19").Value = 13352
$ws.Range("C20").Value = 1
$ws.Range("C21").Value = 121
$ws.Range("C22").Value = 13462
$ws.Range("C23").Value = 0
$ws.Range("C25").Value = 13314
$ws.Range("C26").Value = 0
$ws.Range("C27").Value = 124
$ws.Range("C28").Value = 13405
$ws.Range("C29").Value = 0
$ws.Range("C30").Value = 122
$ws.Range("C31").Value = 13289

$ws = $wb.Worksheets.Item("Selection Sort")
$ws.Range("C4").Value = 4089
$ws.Range("C7").Value = 4070
$ws.Range("C10").Value = 4068
$ws.Range("C13").Value = 4074
$ws.Range("C16").Value = 4063
$ws.Range("C19").Value = 4063
$ws.Range("C21").Value = 41
$ws.Range("C22").Value = 4061
$ws.Range("C24").Value = 41
$ws.Range("C25").Value = 4060
$ws.Range("C28").Value = 4061
$ws.Range("C31").Value = 4058

$ws = $wb.Worksheets.Item("Insertion Sort")
$ws.Range("C4").Value = 1324
$ws.Range("C6").Value = 18
$ws.Range("C7").Value = 1322
$ws.Range("C10").Value = 1315
$ws.Range("C13").Value = 1316
$ws.Range("C19").Value = 1322
$ws.Range("C22").Value = 1323
$ws.Range("C25").Value = 1318
$ws.Range("C28").Value = 1316
$ws.Range("C31").Value = 1327

$ws = $wb.Worksheets.Item("Quick Sort")
$ws.Range("C5").Value = 113
$ws.Range("C6").Value = 234
$ws.Range("C9").Value = 9
$ws.Range("C11").Value = 232
$ws.Range("C15").Value = 114
$ws.Range("C16").Value = 236
$ws.Range("C19").Value = 9
$ws.Range("C20").Value = 114
$ws.Range("C21").Value = 233
$ws.Range("C25").Value = 113
$ws.Range("C26").Value = 234
$ws.Range("C29").Value = 10
$ws.Range("C30").Value = 113
$ws.Range("C31").Value = 234
$ws.Range("C34").Value = 10
$ws.Range("C35").Value = 113
$ws.Range("C36").Value = 234
$ws.Range("C39").Value = 9
$ws.Range("C41").Value = 234
$ws.Range("C45").Value = 113
$ws.Range("C46").Value = 230
$ws.Range("C50").Value = 114
$ws.Range("C51").Value = 232

$ws = $wb.Worksheets.Item("Merge Sort")
$ws.Range("C4").Value = 21
$ws.Range("C5").Value = 194
$ws.Range("C6").Value = 432
$ws.Range("C10").Value = 201
$ws.Range("C11").Value = 421
$ws.Range("C14").Value = 14
$ws.Range("C15").Value = 175
$ws.Range("C16").Value = 453
$ws.Range("C19").Value = 21
$ws.Range("C20").Value = 204
$ws.Range("C25").Value = 180
$ws.Range("C26").Value = 358
$ws.Range("C30").Value = 185
$ws.Range("C31").Value = 511
$ws.Range("C35").Value = 173
$ws.Range("C41").Value = 369
$ws.Range("C45").Value = 173
$ws.Range("C46").Value = 358
$ws.Range("C49").Value = 15
$ws.Range("C50").Value = 174
$ws.Range("C51").Value = 447
